$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-22 12:41:03"

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
